# Commit: "Added 1.1.0 of term"
#
# The "Metadata" sheet holds a simple Property/Value table. Bump the
# "Version" row's value from 1.0.0 -> 1.1.0 and refresh the "Date" row's
# timestamp to match the new release.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$used = $ws.UsedRange
$rows = $used.Rows.Count

for ($r = 1; $r -le $rows; $r++) {
    $label = $ws.Cells.Item($r, 1).Value2

    if ($label -eq "Version") {
        $ws.Cells.Item($r, 2).Value = "1.1.0"
    }
    elseif ($label -eq "Date") {
        $ws.Cells.Item($r, 2).Value = "2023-07-10T23:08:03+02:00"
    }
}
